# Apply "fixed Test Case 1 description" edit to Requirments_Matrix.xlsx (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-23: fill Result / Code Fixed / Regression Result columns
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 5).Value = "pass"   # column E - Result
    $ws.Cells.Item($r, 8).Value = "N/A"    # column H - Code Fixed
    $ws.Cells.Item($r, 9).Value = "pass"   # column I - Regression Result
}

# Row 24 (Test Case 1 description fix): clear the "Test Comments" entry
$ws.Cells.Item(24, 6).Value = ""

# Row 33: replace the blank placeholder comment with the real one
$ws.Cells.Item(33, 10).Value = "File now saves correctly"

# Row 26: add a Regression Test Comment
$ws.Cells.Item(26, 10).Value = "Blank line is now ignored"

# Update the selected/active cell shown when the sheet was last saved
$ws.Activate()
$ws.Range("J26").Select()
